$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (dnn_n51_transpiled.qasm)
$ws.Range("B5").Value = 0.009594336034466857
$ws.Range("C5").Value = 0.004285694231455655
$ws.Range("D5").Value = 242
$ws.Range("E5").Value = 191

# Row 7 (sqrt18.qasm)
$ws.Range("B7").Value = 0.00005818923471113368
$ws.Range("C7").Value = 0.00004027021161707994
$ws.Range("D7").Value = 847
$ws.Range("E7").Value = 781

# Row 8 (dnn_n33_transpiled.qasm)
$ws.Range("B8").Value = 0.1021858996297045
$ws.Range("C8").Value = 0.06999898112563477
$ws.Range("D8").Value = 159
$ws.Range("E8").Value = 122

# Row 9 (qft_n18.qasm)
$ws.Range("B9").Value = 0.03647555862825082
$ws.Range("C9").Value = 0.03172848980237667
$ws.Range("D9").Value = 289
$ws.Range("E9").Value = 264

# Row 10 (DNN16.qasm)
$ws.Range("B10").Value = 0.3970927852420563
$ws.Range("C10").Value = 0.3077348137538805
$ws.Range("D10").Value = 85
$ws.Range("G10").Value = 124

# Row 11 (QV_32.qasm)
$ws.Range("B11").Value = 0.000000004352353397380385
$ws.Range("C11").Value = 0.000000002165762026323673
$ws.Range("D11").Value = 1114
$ws.Range("E11").Value = 1043
$ws.Range("F11").Value = 1486

# Row 13 (hhl_n7.qasm)
$ws.Range("B13").Value = 0.507471706605126
$ws.Range("C13").Value = 0.5004452583181059
$ws.Range("D13").Value = 90
$ws.Range("E13").Value = 81

# Row 15 (google_advantage.qasm)
$ws.Range("B15").Value = 0.7894068617535624
